$d = $word.ActiveDocument

# 1) Merge "Buenos Aires sediments contained " + "huge " into a single run (text unchanged).
$d.Content.Find.Execute("Buenos Aires sediments contained huge ", $true, $false, $false, $false, $false, $true, 1, $false, "Buenos Aires sediments contained huge ", 2)

# 2) Merge " " + "level" + "s" + "," + " similar to sewage sludge" into a single run (text unchanged).
$d.Content.Find.Execute(" levels, similar to sewage sludge", $true, $false, $false, $false, $false, $true, 1, $false, " levels, similar to sewage sludge", 2)

# 3) Change "Higher particle flux and anoxic sediments at Buenos Aires favored sterol preservation "
#    to "Higher fluxes and anoxic sediments at Buenos Aires favored sterol preservation " (merges the
#    whole paragraph into a single run as a side effect).
$d.Content.Find.Execute("Higher particle flux and anoxic sediments at Buenos Aires favored sterol preservation ", $true, $false, $false, $false, $false, $true, 1, $false, "Higher fluxes and anoxic sediments at Buenos Aires favored sterol preservation ", 2)

# 4) Split that run into "Higher flux" | "es" | " and anoxic sediments..." and relocate the _GoBack
#    bookmark between "es" and " and anoxic...".
$rFind = $d.Content
$rFind.Find.Execute("Higher flux", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $rFind.End
$d.Bookmarks.Add("TempSplit", $d.Range($splitPoint, $splitPoint))
$d.Bookmarks("TempSplit").Delete()

$rFind2 = $d.Content
$rFind2.Find.Execute("Higher fluxes", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkPoint = $rFind2.End
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPoint, $bookmarkPoint))
